$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

# Sheet 1
$origStyle = $ws1.Range("B2").Style
$ws1.Range("B2").NumberFormat = "@"
$ws1.Range("B2").Value = "2024.02.14"
$ws1.Range("B2").Style = $origStyle
$origStyle = $ws1.Range("G2").Style
$ws1.Range("G2").NumberFormat = "@"
$ws1.Range("G2").Value = "已停售"
$ws1.Range("G2").Style = $origStyle
$origStyle = $ws1.Range("B3").Style
$ws1.Range("B3").NumberFormat = "@"
$ws1.Range("B3").Value = "2024.02.14"
$ws1.Range("B3").Style = $origStyle
$ws1.Range("F3").Value = 123
$origStyle = $ws1.Range("B4").Style
$ws1.Range("B4").NumberFormat = "@"
$ws1.Range("B4").Value = "2024.02.16"
$ws1.Range("B4").Style = $origStyle
$ws1.Range("F4").Value = 411
$origStyle = $ws1.Range("B5").Style
$ws1.Range("B5").NumberFormat = "@"
$ws1.Range("B5").Value = "2024.02.16"
$ws1.Range("B5").Style = $origStyle
$ws1.Range("F5").Value = 1015
$origStyle = $ws1.Range("B6").Style
$ws1.Range("B6").NumberFormat = "@"
$ws1.Range("B6").Value = "2024.02.17"
$ws1.Range("B6").Style = $origStyle
$ws1.Range("F6").Value = 5698
$origStyle = $ws1.Range("B7").Style
$ws1.Range("B7").NumberFormat = "@"
$ws1.Range("B7").Value = "2024.02.17"
$ws1.Range("B7").Style = $origStyle
$ws1.Range("F7").Value = 528
$origStyle = $ws1.Range("B8").Style
$ws1.Range("B8").NumberFormat = "@"
$ws1.Range("B8").Value = "2024.02.23"
$ws1.Range("B8").Style = $origStyle
$ws1.Range("F8").Value = 755
$origStyle = $ws1.Range("B9").Style
$ws1.Range("B9").NumberFormat = "@"
$ws1.Range("B9").Value = "2024.02.24"
$ws1.Range("B9").Style = $origStyle
$ws1.Range("F9").Value = 986
$origStyle = $ws1.Range("B10").Style
$ws1.Range("B10").NumberFormat = "@"
$ws1.Range("B10").Value = "2024.02.24"
$ws1.Range("B10").Style = $origStyle
$ws1.Range("F10").Value = 831
$origStyle = $ws1.Range("B11").Style
$ws1.Range("B11").NumberFormat = "@"
$ws1.Range("B11").Value = "2024.03.02"
$ws1.Range("B11").Style = $origStyle
$origStyle = $ws1.Range("B12").Style
$ws1.Range("B12").NumberFormat = "@"
$ws1.Range("B12").Value = "2024.03.02"
$ws1.Range("B12").Style = $origStyle
$ws1.Range("F12").Value = 41
$origStyle = $ws1.Range("B13").Style
$ws1.Range("B13").NumberFormat = "@"
$ws1.Range("B13").Value = "2024.03.02"
$ws1.Range("B13").Style = $origStyle
$origStyle = $ws1.Range("B14").Style
$ws1.Range("B14").NumberFormat = "@"
$ws1.Range("B14").Value = "2024.03.02"
$ws1.Range("B14").Style = $origStyle
$origStyle = $ws1.Range("B15").Style
$ws1.Range("B15").NumberFormat = "@"
$ws1.Range("B15").Value = "2024.03.03"
$ws1.Range("B15").Style = $origStyle
$origStyle = $ws1.Range("B16").Style
$ws1.Range("B16").NumberFormat = "@"
$ws1.Range("B16").Value = "2024.03.09"
$ws1.Range("B16").Style = $origStyle
$origStyle = $ws1.Range("B17").Style
$ws1.Range("B17").NumberFormat = "@"
$ws1.Range("B17").Value = "2024.03.09"
$ws1.Range("B17").Style = $origStyle
$ws1.Range("F17").Value = 1925
$origStyle = $ws1.Range("B18").Style
$ws1.Range("B18").NumberFormat = "@"
$ws1.Range("B18").Value = "2024.03.09"
$ws1.Range("B18").Style = $origStyle
$ws1.Range("F18").Value = 1496
$origStyle = $ws1.Range("B19").Style
$ws1.Range("B19").NumberFormat = "@"
$ws1.Range("B19").Value = "2024.03.16"
$ws1.Range("B19").Style = $origStyle
$ws1.Range("F19").Value = 996
$origStyle = $ws1.Range("B20").Style
$ws1.Range("B20").NumberFormat = "@"
$ws1.Range("B20").Value = "2024.03.16"
$ws1.Range("B20").Style = $origStyle
$origStyle = $ws1.Range("B21").Style
$ws1.Range("B21").NumberFormat = "@"
$ws1.Range("B21").Value = "2024.03.16"
$ws1.Range("B21").Style = $origStyle
$ws1.Range("F21").Value = 202
$origStyle = $ws1.Range("B22").Style
$ws1.Range("B22").NumberFormat = "@"
$ws1.Range("B22").Value = "2024.03.16"
$ws1.Range("B22").Style = $origStyle
$ws1.Range("F22").Value = 362
$origStyle = $ws1.Range("B23").Style
$ws1.Range("B23").NumberFormat = "@"
$ws1.Range("B23").Value = "2024.03.16"
$ws1.Range("B23").Style = $origStyle
$ws1.Range("F23").Value = 588
$origStyle = $ws1.Range("B24").Style
$ws1.Range("B24").NumberFormat = "@"
$ws1.Range("B24").Value = "2024.03.16"
$ws1.Range("B24").Style = $origStyle
$ws1.Range("F24").Value = 191
$origStyle = $ws1.Range("B25").Style
$ws1.Range("B25").NumberFormat = "@"
$ws1.Range("B25").Value = "2024.03.16"
$ws1.Range("B25").Style = $origStyle
$ws1.Range("F25").Value = 1063
$origStyle = $ws1.Range("B26").Style
$ws1.Range("B26").NumberFormat = "@"
$ws1.Range("B26").Value = "2024.03.16"
$ws1.Range("B26").Style = $origStyle
$origStyle = $ws1.Range("B27").Style
$ws1.Range("B27").NumberFormat = "@"
$ws1.Range("B27").Value = "2024.03.17"
$ws1.Range("B27").Style = $origStyle
$ws1.Range("F27").Value = 525
$origStyle = $ws1.Range("B28").Style
$ws1.Range("B28").NumberFormat = "@"
$ws1.Range("B28").Value = "2024.03.23"
$ws1.Range("B28").Style = $origStyle
$ws1.Range("F28").Value = 3207
$origStyle = $ws1.Range("B29").Style
$ws1.Range("B29").NumberFormat = "@"
$ws1.Range("B29").Value = "2024.03.23"
$ws1.Range("B29").Style = $origStyle
$ws1.Range("F29").Value = 187
$origStyle = $ws1.Range("B30").Style
$ws1.Range("B30").NumberFormat = "@"
$ws1.Range("B30").Value = "2024.03.24"
$ws1.Range("B30").Style = $origStyle
$ws1.Range("F30").Value = 118
$origStyle = $ws1.Range("B31").Style
$ws1.Range("B31").NumberFormat = "@"
$ws1.Range("B31").Value = "2024.03.24"
$ws1.Range("B31").Style = $origStyle
$ws1.Range("F31").Value = 83
$origStyle = $ws1.Range("B32").Style
$ws1.Range("B32").NumberFormat = "@"
$ws1.Range("B32").Value = "2024.03.24"
$ws1.Range("B32").Style = $origStyle
$origStyle = $ws1.Range("B33").Style
$ws1.Range("B33").NumberFormat = "@"
$ws1.Range("B33").Value = "2024.03.30"
$ws1.Range("B33").Style = $origStyle
$origStyle = $ws1.Range("B34").Style
$ws1.Range("B34").NumberFormat = "@"
$ws1.Range("B34").Value = "2024.03.30"
$ws1.Range("B34").Style = $origStyle
$ws1.Range("F34").Value = 439
$origStyle = $ws1.Range("B35").Style
$ws1.Range("B35").NumberFormat = "@"
$ws1.Range("B35").Value = "2024.03.30"
$ws1.Range("B35").Style = $origStyle
$origStyle = $ws1.Range("B36").Style
$ws1.Range("B36").NumberFormat = "@"
$ws1.Range("B36").Value = "2024.04.04"
$ws1.Range("B36").Style = $origStyle
$origStyle = $ws1.Range("B37").Style
$ws1.Range("B37").NumberFormat = "@"
$ws1.Range("B37").Value = "2024.04.04"
$ws1.Range("B37").Style = $origStyle
$ws1.Range("F37").Value = 17
$origStyle = $ws1.Range("B38").Style
$ws1.Range("B38").NumberFormat = "@"
$ws1.Range("B38").Value = "2024.04.04"
$ws1.Range("B38").Style = $origStyle
$origStyle = $ws1.Range("B39").Style
$ws1.Range("B39").NumberFormat = "@"
$ws1.Range("B39").Value = "2024.04.05"
$ws1.Range("B39").Style = $origStyle
$origStyle = $ws1.Range("B40").Style
$ws1.Range("B40").NumberFormat = "@"
$ws1.Range("B40").Value = "2024.04.05"
$ws1.Range("B40").Style = $origStyle
$ws1.Range("F40").Value = 769
$origStyle = $ws1.Range("B41").Style
$ws1.Range("B41").NumberFormat = "@"
$ws1.Range("B41").Value = "2024.04.05"
$ws1.Range("B41").Style = $origStyle
$origStyle = $ws1.Range("B42").Style
$ws1.Range("B42").NumberFormat = "@"
$ws1.Range("B42").Value = "2024.04.13"
$ws1.Range("B42").Style = $origStyle
$origStyle = $ws1.Range("B43").Style
$ws1.Range("B43").NumberFormat = "@"
$ws1.Range("B43").Value = "2024.04.20"
$ws1.Range("B43").Style = $origStyle
$ws1.Range("F43").Value = 68
$origStyle = $ws1.Range("B44").Style
$ws1.Range("B44").NumberFormat = "@"
$ws1.Range("B44").Value = "2024.04.20"
$ws1.Range("B44").Style = $origStyle

# Sheet 2
$origStyle = $ws2.Range("B2").Style
$ws2.Range("B2").NumberFormat = "@"
$ws2.Range("B2").Value = "2024.02.17"
$ws2.Range("B2").Style = $origStyle
$origStyle = $ws2.Range("B3").Style
$ws2.Range("B3").NumberFormat = "@"
$ws2.Range("B3").Value = "2024.02.24"
$ws2.Range("B3").Style = $origStyle
$origStyle = $ws2.Range("B4").Style
$ws2.Range("B4").NumberFormat = "@"
$ws2.Range("B4").Value = "2024.02.24"
$ws2.Range("B4").Style = $origStyle
$ws2.Range("F4").Value = 232
$origStyle = $ws2.Range("B5").Style
$ws2.Range("B5").NumberFormat = "@"
$ws2.Range("B5").Value = "2024.02.25"
$ws2.Range("B5").Style = $origStyle
$origStyle = $ws2.Range("B6").Style
$ws2.Range("B6").NumberFormat = "@"
$ws2.Range("B6").Value = "2024.02.25"
$ws2.Range("B6").Style = $origStyle
$ws2.Range("F6").Value = 157
$origStyle = $ws2.Range("B7").Style
$ws2.Range("B7").NumberFormat = "@"
$ws2.Range("B7").Value = "2024.03.08"
$ws2.Range("B7").Style = $origStyle
$origStyle = $ws2.Range("B8").Style
$ws2.Range("B8").NumberFormat = "@"
$ws2.Range("B8").Value = "2024.03.16"
$ws2.Range("B8").Style = $origStyle
$origStyle = $ws2.Range("B9").Style
$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "2024.04.21"
$ws2.Range("B9").Style = $origStyle
$origStyle = $ws2.Range("B10").Style
$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "2024.04.27"
$ws2.Range("B10").Style = $origStyle
$origStyle = $ws2.Range("B11").Style
$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "2024.05.01"
$ws2.Range("B11").Style = $origStyle

# Sheet 3
$origStyle = $ws3.Range("B2").Style
$ws3.Range("B2").NumberFormat = "@"
$ws3.Range("B2").Value = "2023.12.30"
$ws3.Range("B2").Style = $origStyle

# Sheet 4
$origStyle = $ws4.Range("B2").Style
$ws4.Range("B2").NumberFormat = "@"
$ws4.Range("B2").Value = "2023.12.30"
$ws4.Range("B2").Style = $origStyle
$origStyle = $ws4.Range("B3").Style
$ws4.Range("B3").NumberFormat = "@"
$ws4.Range("B3").Value = "2024.02.14"
$ws4.Range("B3").Style = $origStyle
$ws4.Range("F3").Value = 123
$origStyle = $ws4.Range("B4").Style
$ws4.Range("B4").NumberFormat = "@"
$ws4.Range("B4").Value = "2024.02.16"
$ws4.Range("B4").Style = $origStyle
$ws4.Range("F4").Value = 1015
$origStyle = $ws4.Range("B5").Style
$ws4.Range("B5").NumberFormat = "@"
$ws4.Range("B5").Value = "2024.02.17"
$ws4.Range("B5").Style = $origStyle
$origStyle = $ws4.Range("B6").Style
$ws4.Range("B6").NumberFormat = "@"
$ws4.Range("B6").Value = "2024.02.17"
$ws4.Range("B6").Style = $origStyle
$ws4.Range("F6").Value = 5698
$origStyle = $ws4.Range("B7").Style
$ws4.Range("B7").NumberFormat = "@"
$ws4.Range("B7").Value = "2024.02.17"
$ws4.Range("B7").Style = $origStyle
$ws4.Range("F7").Value = 528
$origStyle = $ws4.Range("B8").Style
$ws4.Range("B8").NumberFormat = "@"
$ws4.Range("B8").Value = "2024.02.23"
$ws4.Range("B8").Style = $origStyle
$ws4.Range("F8").Value = 755
$origStyle = $ws4.Range("B9").Style
$ws4.Range("B9").NumberFormat = "@"
$ws4.Range("B9").Value = "2024.02.24"
$ws4.Range("B9").Style = $origStyle
$origStyle = $ws4.Range("B10").Style
$ws4.Range("B10").NumberFormat = "@"
$ws4.Range("B10").Value = "2024.02.24"
$ws4.Range("B10").Style = $origStyle
$ws4.Range("F10").Value = 232
$origStyle = $ws4.Range("B11").Style
$ws4.Range("B11").NumberFormat = "@"
$ws4.Range("B11").Value = "2024.02.24"
$ws4.Range("B11").Style = $origStyle
$ws4.Range("F11").Value = 986
$origStyle = $ws4.Range("B12").Style
$ws4.Range("B12").NumberFormat = "@"
$ws4.Range("B12").Value = "2024.02.24"
$ws4.Range("B12").Style = $origStyle
$ws4.Range("F12").Value = 831
$origStyle = $ws4.Range("B13").Style
$ws4.Range("B13").NumberFormat = "@"
$ws4.Range("B13").Value = "2024.02.25"
$ws4.Range("B13").Style = $origStyle
$origStyle = $ws4.Range("B14").Style
$ws4.Range("B14").NumberFormat = "@"
$ws4.Range("B14").Value = "2024.02.25"
$ws4.Range("B14").Style = $origStyle
$ws4.Range("F14").Value = 157
$origStyle = $ws4.Range("B15").Style
$ws4.Range("B15").NumberFormat = "@"
$ws4.Range("B15").Value = "2024.03.02"
$ws4.Range("B15").Style = $origStyle
$origStyle = $ws4.Range("B16").Style
$ws4.Range("B16").NumberFormat = "@"
$ws4.Range("B16").Value = "2024.03.02"
$ws4.Range("B16").Style = $origStyle
$ws4.Range("F16").Value = 41
$origStyle = $ws4.Range("B17").Style
$ws4.Range("B17").NumberFormat = "@"
$ws4.Range("B17").Value = "2024.03.02"
$ws4.Range("B17").Style = $origStyle
$origStyle = $ws4.Range("B18").Style
$ws4.Range("B18").NumberFormat = "@"
$ws4.Range("B18").Value = "2024.03.02"
$ws4.Range("B18").Style = $origStyle
$origStyle = $ws4.Range("B19").Style
$ws4.Range("B19").NumberFormat = "@"
$ws4.Range("B19").Value = "2024.03.03"
$ws4.Range("B19").Style = $origStyle
$origStyle = $ws4.Range("B20").Style
$ws4.Range("B20").NumberFormat = "@"
$ws4.Range("B20").Value = "2024.03.08"
$ws4.Range("B20").Style = $origStyle
$origStyle = $ws4.Range("B21").Style
$ws4.Range("B21").NumberFormat = "@"
$ws4.Range("B21").Value = "2024.03.09"
$ws4.Range("B21").Style = $origStyle
$origStyle = $ws4.Range("B22").Style
$ws4.Range("B22").NumberFormat = "@"
$ws4.Range("B22").Value = "2024.03.09"
$ws4.Range("B22").Style = $origStyle
$ws4.Range("F22").Value = 1925
$origStyle = $ws4.Range("B23").Style
$ws4.Range("B23").NumberFormat = "@"
$ws4.Range("B23").Value = "2024.03.09"
$ws4.Range("B23").Style = $origStyle
$ws4.Range("F23").Value = 1496
$origStyle = $ws4.Range("B24").Style
$ws4.Range("B24").NumberFormat = "@"
$ws4.Range("B24").Value = "2024.03.16"
$ws4.Range("B24").Style = $origStyle
$ws4.Range("F24").Value = 996
$origStyle = $ws4.Range("B25").Style
$ws4.Range("B25").NumberFormat = "@"
$ws4.Range("B25").Value = "2024.03.16"
$ws4.Range("B25").Style = $origStyle
$ws4.Range("F25").Value = 202
$origStyle = $ws4.Range("B26").Style
$ws4.Range("B26").NumberFormat = "@"
$ws4.Range("B26").Value = "2024.03.16"
$ws4.Range("B26").Style = $origStyle
$ws4.Range("F26").Value = 362
$origStyle = $ws4.Range("B27").Style
$ws4.Range("B27").NumberFormat = "@"
$ws4.Range("B27").Value = "2024.03.16"
$ws4.Range("B27").Style = $origStyle
$origStyle = $ws4.Range("B28").Style
$ws4.Range("B28").NumberFormat = "@"
$ws4.Range("B28").Value = "2024.03.16"
$ws4.Range("B28").Style = $origStyle
$ws4.Range("F28").Value = 588
$origStyle = $ws4.Range("B29").Style
$ws4.Range("B29").NumberFormat = "@"
$ws4.Range("B29").Value = "2024.03.16"
$ws4.Range("B29").Style = $origStyle
$ws4.Range("F29").Value = 191
$origStyle = $ws4.Range("B30").Style
$ws4.Range("B30").NumberFormat = "@"
$ws4.Range("B30").Value = "2024.03.16"
$ws4.Range("B30").Style = $origStyle
$ws4.Range("F30").Value = 1063
$origStyle = $ws4.Range("B31").Style
$ws4.Range("B31").NumberFormat = "@"
$ws4.Range("B31").Value = "2024.03.23"
$ws4.Range("B31").Style = $origStyle
$ws4.Range("F31").Value = 3207
$origStyle = $ws4.Range("B32").Style
$ws4.Range("B32").NumberFormat = "@"
$ws4.Range("B32").Value = "2024.03.23"
$ws4.Range("B32").Style = $origStyle
$ws4.Range("F32").Value = 187
$origStyle = $ws4.Range("B33").Style
$ws4.Range("B33").NumberFormat = "@"
$ws4.Range("B33").Value = "2024.03.24"
$ws4.Range("B33").Style = $origStyle
$ws4.Range("F33").Value = 118
$origStyle = $ws4.Range("B34").Style
$ws4.Range("B34").NumberFormat = "@"
$ws4.Range("B34").Value = "2024.03.24"
$ws4.Range("B34").Style = $origStyle
$ws4.Range("F34").Value = 83
$origStyle = $ws4.Range("B35").Style
$ws4.Range("B35").NumberFormat = "@"
$ws4.Range("B35").Value = "2024.03.24"
$ws4.Range("B35").Style = $origStyle
$origStyle = $ws4.Range("B36").Style
$ws4.Range("B36").NumberFormat = "@"
$ws4.Range("B36").Value = "2024.03.30"
$ws4.Range("B36").Style = $origStyle
$origStyle = $ws4.Range("B37").Style
$ws4.Range("B37").NumberFormat = "@"
$ws4.Range("B37").Value = "2024.03.30"
$ws4.Range("B37").Style = $origStyle
$ws4.Range("F37").Value = 439
$origStyle = $ws4.Range("B38").Style
$ws4.Range("B38").NumberFormat = "@"
$ws4.Range("B38").Value = "2024.03.30"
$ws4.Range("B38").Style = $origStyle
$origStyle = $ws4.Range("B39").Style
$ws4.Range("B39").NumberFormat = "@"
$ws4.Range("B39").Value = "2024.04.04"
$ws4.Range("B39").Style = $origStyle
$origStyle = $ws4.Range("B40").Style
$ws4.Range("B40").NumberFormat = "@"
$ws4.Range("B40").Value = "2024.04.04"
$ws4.Range("B40").Style = $origStyle
$ws4.Range("F40").Value = 17
$origStyle = $ws4.Range("B41").Style
$ws4.Range("B41").NumberFormat = "@"
$ws4.Range("B41").Value = "2024.04.05"
$ws4.Range("B41").Style = $origStyle
$origStyle = $ws4.Range("B42").Style
$ws4.Range("B42").NumberFormat = "@"
$ws4.Range("B42").Value = "2024.04.05"
$ws4.Range("B42").Style = $origStyle
$ws4.Range("F42").Value = 769
$origStyle = $ws4.Range("B43").Style
$ws4.Range("B43").NumberFormat = "@"
$ws4.Range("B43").Value = "2024.04.05"
$ws4.Range("B43").Style = $origStyle
$origStyle = $ws4.Range("B44").Style
$ws4.Range("B44").NumberFormat = "@"
$ws4.Range("B44").Value = "2024.04.13"
$ws4.Range("B44").Style = $origStyle
$origStyle = $ws4.Range("B45").Style
$ws4.Range("B45").NumberFormat = "@"
$ws4.Range("B45").Value = "2024.04.20"
$ws4.Range("B45").Style = $origStyle
$ws4.Range("F45").Value = 68
$origStyle = $ws4.Range("B46").Style
$ws4.Range("B46").NumberFormat = "@"
$ws4.Range("B46").Value = "2024.04.20"
$ws4.Range("B46").Style = $origStyle
$origStyle = $ws4.Range("B47").Style
$ws4.Range("B47").NumberFormat = "@"
$ws4.Range("B47").Value = "2024.04.21"
$ws4.Range("B47").Style = $origStyle
$origStyle = $ws4.Range("B48").Style
$ws4.Range("B48").NumberFormat = "@"
$ws4.Range("B48").Value = "2024.04.27"
$ws4.Range("B48").Style = $origStyle
$origStyle = $ws4.Range("B49").Style
$ws4.Range("B49").NumberFormat = "@"
$ws4.Range("B49").Value = "2024.05.01"
$ws4.Range("B49").Style = $origStyle
